$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.289.07'
$ws.Range('E2').Value = '  +0.25%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.870.54'
$ws.Range('E3').Value = '  +0.59%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '234.97'
$ws.Range('E5').Value = '  -0.55%  '
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4698'
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2869'
$ws.Range('E8').Value = '  +0.19%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06590'
$ws.Range('E9').Value = '  +0.73%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.61'
$ws.Range('E10').Value = '  -0.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07952'
$ws.Range('E11').Value = '  +0.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '96.41'
$ws.Range('E12').Value = '  -0.80%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.873.97'
$ws.Range('E13').Value = '  +0.71%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6956'
$ws.Range('E14').Value = '  +2.15%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.106'
$ws.Range('E15').Value = '  -1.08%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '268.39'
$ws.Range('E16').Value = '  -0.48%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.328.05'
$ws.Range('E17').Value = '  +0.41%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.10'
$ws.Range('E18').Value = '  +4.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007760'
$ws.Range('E19').Value = '  +5.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('E20').Value = '  +0.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.119.86'
$ws.Range('E21').Value = '  +0.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.002'
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.267'
$ws.Range('E23').Value = '  -0.91%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.206'
$ws.Range('E24').Value = '  +0.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.372'
$ws.Range('E25').Value = '  +1.60%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '167.57'
$ws.Range('E26').Value = '  +0.39%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.85'
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.950'
$ws.Range('E28').Value = '  -0.64%  '
$ws.Range('E29').Value = '  -1.26%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09888'
$ws.Range('E30').Value = '  +0.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.332'
$ws.Range('E31').Value = '  -0.65%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.462'
$ws.Range('E32').Value = '  -1.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.066'
$ws.Range('E33').Value = '  +0.43%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04712'
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.134'
$ws.Range('E35').Value = '  +0.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7010'
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.726'
$ws.Range('E37').Value = '  +0.61%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01870'
$ws.Range('E38').Value = '  -0.43%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.794'
$ws.Range('E39').Value = '  +6.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.234'
$ws.Range('E40').Value = '  -0.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '71.76'
$ws.Range('E41').Value = '  -5.20%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.953'
$ws.Range('E42').Value = '  +0.47%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4175'
$ws.Range('E43').Value = '  +0.36%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.000'
$ws.Range('E44').Value = '  -0.03%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8401'
$ws.Range('E45').Value = '  -1.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '102.57'
$ws.Range('E46').Value = '  -0.82%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.103'
$ws.Range('E47').Value = '  -0.91%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.145'
$ws.Range('E48').Value = '  -1.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '915.81'
$ws.Range('E49').Value = '  -4.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.51'
$ws.Range('E50').Value = '  +1.17%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05681'
$ws.Range('E51').Value = '  +0.51%  '
